$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new value, and whether the value needs a
# leading quote-prefix to stop Excel from auto-parsing a numeric-
# looking string (e.g. "1.000" -> 1) when typed into the cell.
$updates = @(
    @{ Cell = 'D2'; Value = '30.208.44'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +3.31%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.895.06'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  -0.36%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '325.20'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +3.39%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.5162'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  +0.38%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.4009'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  +2.48%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.08444'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  +0.43%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '42.66'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +0.99%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '1.116'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  +0.27%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '23.25'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  +12.74%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '6.424'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  +2.81%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '1.892.82'; ForceText = $false },
    @{ Cell = 'E14'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '7.335'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  +0.38%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -0.33%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '94.69'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  +1.76%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '0.00001111'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  +0.53%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.06649'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -1.40%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '18.23'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +2.35%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '1.000'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -0.35%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '5.945'; ForceText = $true },
    @{ Cell = 'D23'; Value = '30.211.87'; ForceText = $false },
    @{ Cell = 'E23'; Value = '  +3.29%  '; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +1.61%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '2.229'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +0.51%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.108.87'; ForceText = $false },
    @{ Cell = 'E26'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '21.68'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  +4.02%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '161.17'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +1.15%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '2.370'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -2.23%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '128.86'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  +0.92%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '1.096'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  +3.82%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.1057'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  +1.04%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '6.049'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -2.04%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '3.757'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  +2.62%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.02490'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +0.44%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.06554'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -0.17%  '; ForceText = $false },
    @{ Cell = 'B37'; Value = 'InternetComputer(DFINITY)'; ForceText = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false },
    @{ Cell = 'D37'; Value = '5.249'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +1.34%  '; ForceText = $false },
    @{ Cell = 'B38'; Value = 'Algorand'; ForceText = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.2200'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +0.58%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '1.217'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -0.69%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '11.78'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +4.74%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.6499'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  +0.10%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '8.708'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -3.55%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.236'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +0.37%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '0.6104'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  +1.01%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '13.29'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +0.73%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '3.703'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +0.91%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '2.055'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +0.32%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '1.234'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +0.52%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '124.33'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +0.75%  '; ForceText = $false },
    @{ Cell = 'E50'; Value = '  -0.74%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '78.84'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +1.92%  '; ForceText = $false }
)

$textForcedCells = @()

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).Value = "'" + $u.Value
        $textForcedCells += $u.Cell
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

# Reset style back to Normal for the quote-prefixed cells so they keep
# the workbook's original (default) cell style instead of picking up
# a stray quotePrefix-flagged style.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}

